$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the grade values for row 8 (student #1)
$ws.Range("G8").Value = 0.83
$ws.Range("H8").Value = 0.31

# Update the active selection to match the author's final cursor position
$ws.Range("G9").Select()
